# Update cryptos list: price (D) and volume (E) changes, plus a re-sort of rows 43-46
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name / link updates (rows 43-46 were re-ordered) ---
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("B44").Value = "PancakeSwap"
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("C44").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"

# --- Price column (D) updates; force text so values like "20.456.95" or
#     "1.012" are not coerced into numbers/dates by Excel ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.456.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.468.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "280.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.8929"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3718"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3193"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "40.08"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.052"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06645"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.007"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.549"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.217"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.478.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001033"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.05681"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.8967"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.692"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.307"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.735.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.284"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "137.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.639.86"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "113.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.968"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.211"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8458"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07809"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06127"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.477"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.871"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.172"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02049"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1876"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9180"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5370"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.584"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.886"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "123.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5274"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.824"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06442"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.040"
$ws.Range("D51").Style = "Normal"

# --- Volume(1h) column (E) updates ---
$ws.Range("E2").Value = "  +2.71%  "
$ws.Range("E3").Value = "  +4.63%  "
$ws.Range("E4").Value = "  +0.82%  "
$ws.Range("E5").Value = "  +2.95%  "
$ws.Range("E6").Value = "  -11.02%  "
$ws.Range("E7").Value = "  +0.80%  "
$ws.Range("E8").Value = "  +4.38%  "
$ws.Range("E9").Value = "  +3.08%  "
$ws.Range("E10").Value = "  +6.50%  "
$ws.Range("E11").Value = "  +1.88%  "
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("E13").Value = "  +2.97%  "
$ws.Range("E14").Value = "  +6.95%  "
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("E16").Value = "  +4.83%  "
$ws.Range("E17").Value = "  +3.10%  "
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("E19").Value = "  -10.63%  "
$ws.Range("E20").Value = "  -3.38%  "
$ws.Range("E22").Value = "  +2.01%  "
$ws.Range("E23").Value = "  +4.58%  "
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("E25").Value = "  +4.08%  "
$ws.Range("E26").Value = "  +1.84%  "
$ws.Range("E27").Value = "  -0.45%  "
$ws.Range("E28").Value = "  +3.67%  "
$ws.Range("E29").Value = "  +4.42%  "
$ws.Range("E30").Value = "  +4.80%  "
$ws.Range("E31").Value = "  +2.85%  "
$ws.Range("E32").Value = "  -2.41%  "
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("E34").Value = "  +1.68%  "
$ws.Range("E35").Value = "  +6.84%  "
$ws.Range("E36").Value = "  +16.48%  "
$ws.Range("E37").Value = "  +2.12%  "
$ws.Range("E38").Value = "  +11.12%  "
$ws.Range("E39").Value = "  +4.73%  "
$ws.Range("E40").Value = "  +1.58%  "
$ws.Range("E41").Value = "  -1.20%  "
$ws.Range("E42").Value = "  -8.40%  "
$ws.Range("E43").Value = "  +2.37%  "
$ws.Range("E44").Value = "  +1.78%  "
$ws.Range("E45").Value = "  +1.83%  "
$ws.Range("E46").Value = "  -17.93%  "
$ws.Range("E47").Value = "  +13.01%  "
$ws.Range("E48").Value = "  +3.95%  "
$ws.Range("E49").Value = "  +2.24%  "
$ws.Range("E50").Value = "  +5.25%  "
$ws.Range("E51").Value = "  -0.12%  "
